# The "Trigger" / "erledigt" columns of the notification-codes table
# contained standalone "Ja"/"ja"/"Nein"/"nein" cell values that are no
# longer needed - clear them out (leaving the now-empty paragraph/cell
# behind), matching the "refactor and clean up" commit.

$d = $word.ActiveDocument

$values = @("Ja", "ja", "Nein", "nein")

foreach ($value in $values) {
    # MatchCase + MatchWholeWord so we only ever hit the standalone
    # Ja/ja/Nein/nein table-cell contents, never a substring of some
    # other word. Replace = 2 (wdReplaceAll) sweeps every occurrence in
    # $d.Content in one call.
    $d.Content.Find.Execute($value, $true, $true, $false, $false, $false, `
                             $true, 1, $false, "", 2)
}
